$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.901.62"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "2.586.34"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'587.87"
$ws.Range("E5").Value = "  -1.61%  "
$ws.Range("D6").Value = "'148.16"
$ws.Range("E6").Value = "  -4.26%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "2.586.42"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("E10").Value = "  -4.02%  "
$ws.Range("D11").Value = "'0.160"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "'5.13"
$ws.Range("E12").Value = "  -2.33%  "
$ws.Range("E13").Value = "  -3.90%  "
$ws.Range("D14").Value = "'26.93"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").Value = "3.061.51"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = "  -5.31%  "
$ws.Range("D17").Value = "66.774.29"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "2.596.94"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "'362.56"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").Value = "'7.27"
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'71.91"
$ws.Range("E26").Value = "  +6.50%  "
$ws.Range("D27").Value = "'9.83"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").Value = "2.711.62"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'571.08"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").Value = "0.0₃0967"
$ws.Range("E31").Value = "  -6.89%  "
$ws.Range("E32").Value = "  -5.25%  "
$ws.Range("E33").Value = "  -4.63%  "
$ws.Range("D34").Value = "'1.79"
$ws.Range("E34").Value = "  -3.57%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -7.03%  "
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").Value = "'156.27"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("D39").Value = "'18.81"
$ws.Range("E39").Value = "  -2.83%  "
$ws.Range("D40").Value = "'0.363"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("E41").Value = "  -1.96%  "
$ws.Range("D42").Value = "'5.12"
$ws.Range("E42").Value = "  -3.93%  "
$ws.Range("D43").Value = "'16.77"
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("D44").Value = "'2.46"
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("D46").Value = "'152.56"
$ws.Range("D47").Value = "0.0₆0279"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("D50").Value = "'1.67"
$ws.Range("E50").Value = "  -3.24%  "
$ws.Range("D51").Value = "'21.10"
$ws.Range("E51").Value = "  +0.26%  "
